# "Export with no is_pref and no lev distance"
#
# The rows' speaker-variant text (columns B/C) is regenerated (here it
# rotates: row2's old text moves to row5, rows 3/4/5 shift up into
# 2/3/4), and the is_prefered flag (column D, previously "x" on every
# row) is cleared on every row since the export no longer carries
# is_pref / levenshtein-distance info.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @{ B = "#he,-sa,-vo,-been-marketroy"; C = "He, Sa, vo, been marketroy" },
  @{ B = "#inquam-koenat,-rekoolatt,-katan,-extra,-struynalibol"; C = "Inquam koenat, Rekoolatt, katan, Extra, struynalibol" },
  @{ B = "#hier-leit-hy-die-het-spijt-dat-hem-de-doodt-verriedt:"; C = "Hier leit hy die het spijt dat hem de doodt verriedt:" },
  @{ B = "#he,-ho-onsikokx-kalandarina,-le,-bo,-bonsibokx,`nmalandarina-ronsikfokx"; C = "He, ho onsikokx Kalandarina, Le, bo, bonsibokx,`nmalandarina Ronsikfokx" }
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 2).Value = $data[$i].B
  $ws.Cells.Item($row, 3).Value = $data[$i].C

  # Clear the is_prefered cell, but keep it as a (blank) text cell rather
  # than deleting it outright: assigning "" removes the cell entirely in
  # this engine, so assign a lone quote-prefix character (Excel's "force
  # text" marker), which stores as an empty string of type Text, then
  # strip the quote-prefix cell format it introduces so no stray style
  # is left behind.
  $ws.Cells.Item($row, 4).Value = "'"
  $ws.Cells.Item($row, 4).Style = "Normal"
}
